# Scheduled price-refresh update: re-computed profit/price figures
# (currentAveragePrice*, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1927.625
$ws.Range("I6").Value = 636.8333
$ws.Range("K6").Value = 1910.4999
$ws.Range("M6").Value = -1798.4999
$ws.Range("H17").Value = 3417.4443
$ws.Range("J17").Value = 3417.4443
$ws.Range("L17").Value = 10252.3329
$ws.Range("N17").Value = -10588.3329
$ws.Range("H92").Value = 189.78947
$ws.Range("I92").Value = 113.875
$ws.Range("K92").Value = 113.875
$ws.Range("M92").Value = 1134.125
$ws.Range("H96").Value = 633.25
$ws.Range("I96").Value = 273.85715
$ws.Range("K96").Value = 821.5714499999999
$ws.Range("M96").Value = 551.4285500000001
$ws.Range("H100").Value = 2696
$ws.Range("I100").Value = 1544.5
$ws.Range("K100").Value = 1544.5
$ws.Range("M100").Value = -1003.5
$ws.Range("H111").Value = 101854.9
$ws.Range("J111").Value = 334932.66
$ws.Range("L111").Value = 1004797.98
$ws.Range("N111").Value = -1010931.98
$ws.Range("H119").Value = 1204.3077
$ws.Range("J119").Value = 1204.3077
$ws.Range("L119").Value = 3612.9231
$ws.Range("N119").Value = -13288.9231
$ws.Range("H127").Value = 2154.75
$ws.Range("I127").Value = 793.9167
$ws.Range("K127").Value = 2381.7501
$ws.Range("M127").Value = 2578.2499
$ws.Range("H132").Value = 3302.4285
$ws.Range("I132").Value = 3302.4285
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9907.2855
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7377.2855
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 3209.6
$ws.Range("I137").Value = 2658.92
$ws.Range("J137").Value = 4127.4
$ws.Range("K137").Value = 7976.76
$ws.Range("L137").Value = 12382.2
$ws.Range("M137").Value = -5426.76
$ws.Range("N137").Value = -17482.2
$ws.Range("H138").Value = 4513.457
$ws.Range("I138").Value = 4814.636
$ws.Range("J138").Value = 4401.1523
$ws.Range("K138").Value = 14443.908
$ws.Range("L138").Value = 13203.4569
$ws.Range("M138").Value = -9303.908000000001
$ws.Range("N138").Value = -23483.4569

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 416.6111
$ws.Range("I5").Value = 411.7647
$ws.Range("K5").Value = 411.7647
$ws.Range("M5").Value = -299.7647
$ws.Range("H32").Value = 16296.9
$ws.Range("I32").Value = 11634.283
$ws.Range("J32").Value = 51599.57
$ws.Range("K32").Value = 11634.283
$ws.Range("L32").Value = 51599.57
$ws.Range("M32").Value = -11347.283
$ws.Range("N32").Value = -52173.57
$ws.Range("H56").Value = 6933.3335
$ws.Range("I56").Value = 6933.3335
$ws.Range("K56").Value = 6933.3335
$ws.Range("M56").Value = -6191.3335
$ws.Range("H74").Value = 6705.25
$ws.Range("I74").Value = 3611.7058
$ws.Range("K74").Value = 3611.7058
$ws.Range("M74").Value = -2737.7058
$ws.Range("H77").Value = 6705.25
$ws.Range("I77").Value = 3611.7058
$ws.Range("K77").Value = 18058.529
$ws.Range("M77").Value = -13690.529
$ws.Range("H102").Value = 1100
$ws.Range("I102").Value = 1100
$ws.Range("K102").Value = 1100
$ws.Range("M102").Value = 522
$ws.Range("H132").Value = 6549.591
$ws.Range("I132").Value = 2412.9375
$ws.Range("K132").Value = 7238.8125
$ws.Range("M132").Value = -4708.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 416.6111
$ws.Range("I4").Value = 411.7647
$ws.Range("K4").Value = 411.7647
$ws.Range("M4").Value = -296.7647
$ws.Range("H86").Value = 2569.6667
$ws.Range("I86").Value = 2569.6667
$ws.Range("K86").Value = 2569.6667
$ws.Range("M86").Value = -1446.6667
$ws.Range("H89").Value = 2569.6667
$ws.Range("I89").Value = 2569.6667
$ws.Range("K89").Value = 12848.3335
$ws.Range("M89").Value = -7232.333500000001
$ws.Range("H134").Value = 4868.321
$ws.Range("I134").Value = 3991.9546
$ws.Range("K134").Value = 11975.8638
$ws.Range("M134").Value = -9440.863799999999
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2060.5715
$ws.Range("I22").Value = 1100.25
$ws.Range("J22").Value = 2444.7
$ws.Range("K22").Value = 1100.25
$ws.Range("L22").Value = 2444.7
$ws.Range("M22").Value = -750.25
$ws.Range("N22").Value = -3144.7
$ws.Range("H31").Value = 8948.579
$ws.Range("I31").Value = 4786.55
$ws.Range("J31").Value = 13573.056
$ws.Range("K31").Value = 4786.55
$ws.Range("L31").Value = 13573.056
$ws.Range("M31").Value = -4491.55
$ws.Range("N31").Value = -14163.056
$ws.Range("H34").Value = 8948.579
$ws.Range("I34").Value = 4786.55
$ws.Range("J34").Value = 13573.056
$ws.Range("K34").Value = 4786.55
$ws.Range("L34").Value = 13573.056
$ws.Range("M34").Value = -4584.55
$ws.Range("N34").Value = -13977.056
$ws.Range("H105").Value = 1986.9231
$ws.Range("I105").Value = 2069.1667
$ws.Range("K105").Value = 2069.1667
$ws.Range("M105").Value = -322.1667000000002
$ws.Range("H107").Value = 1399.9333
$ws.Range("I107").Value = 1082.4286
$ws.Range("J107").Value = 2140.7778
$ws.Range("K107").Value = 1082.4286
$ws.Range("L107").Value = 2140.7778
$ws.Range("M107").Value = 837.5714
$ws.Range("N107").Value = -5980.7778
$ws.Range("H111").Value = 275000
$ws.Range("J111").Value = 275000
$ws.Range("L111").Value = 275000
$ws.Range("N111").Value = -283180
$ws.Range("H141").Value = 288024.6
$ws.Range("J141").Value = 325635.75
$ws.Range("L141").Value = 325635.75
$ws.Range("N141").Value = -335995.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 3165.6667
$ws.Range("I46").Value = 433.33334
$ws.Range("J46").Value = 4531.8335
$ws.Range("K46").Value = 1300.00002
$ws.Range("L46").Value = 13595.5005
$ws.Range("M46").Value = -1209.00002
$ws.Range("N46").Value = -13777.5005
$ws.Range("H121").Value = 1055475.4
$ws.Range("I121").Value = 215.2
$ws.Range("K121").Value = 645.5999999999999
$ws.Range("M121").Value = 664.4000000000001
$ws.Range("H128").Value = 257280.58
$ws.Range("I128").Value = 257280.58
$ws.Range("K128").Value = 771841.74
$ws.Range("M128").Value = -766861.74

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5249.1113
$ws.Range("I70").Value = 4500
$ws.Range("K70").Value = 4500
$ws.Range("M70").Value = -4230
$ws.Range("H73").Value = 5249.1113
$ws.Range("I73").Value = 4500
$ws.Range("K73").Value = 4500
$ws.Range("M73").Value = -3564
$ws.Range("H80").Value = 1640
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 1780
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 1780
$ws.Range("M80").Value = -502
$ws.Range("N80").Value = -3776
$ws.Range("H83").Value = 1640
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 1780
$ws.Range("K83").Value = 7500
$ws.Range("L83").Value = 8900
$ws.Range("M83").Value = -2508
$ws.Range("N83").Value = -18884
$ws.Range("H97").Value = 923.8125
$ws.Range("I97").Value = 904.1
$ws.Range("J97").Value = 956.6667
$ws.Range("K97").Value = 904.1
$ws.Range("L97").Value = 956.6667
$ws.Range("M97").Value = -408.1
$ws.Range("N97").Value = -1948.6667
$ws.Range("H122").Value = 7208.2607
$ws.Range("I122").Value = 4253.2144
$ws.Range("K122").Value = 12759.6432
$ws.Range("M122").Value = -10309.6432
$ws.Range("H126").Value = 7480.4688
$ws.Range("I126").Value = 6785.95
$ws.Range("J126").Value = 8638
$ws.Range("K126").Value = 20357.85
$ws.Range("L126").Value = 25914
$ws.Range("M126").Value = -17887.85
$ws.Range("N126").Value = -30854
$ws.Range("H132").Value = 5222.396
$ws.Range("I132").Value = 5201.85
$ws.Range("K132").Value = 15605.55
$ws.Range("M132").Value = -13075.55

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 648.1070999999999
$ws.Range("I16").Value = 557.6818
$ws.Range("K16").Value = 557.6818
$ws.Range("M16").Value = -387.6818
$ws.Range("H22").Value = 2366.6667
$ws.Range("I22").Value = 2037.5
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 2037.5
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -1742.5
$ws.Range("N22").Value = -5590
$ws.Range("H27").Value = 2366.6667
$ws.Range("I27").Value = 2037.5
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 2037.5
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -1930.5
$ws.Range("N27").Value = -5214
$ws.Range("H38").Value = 36000
$ws.Range("J38").Value = 36000
$ws.Range("L38").Value = 36000
$ws.Range("N38").Value = -36820
$ws.Range("H40").Value = 4759.28
$ws.Range("I40").Value = 3618.1904
$ws.Range("K40").Value = 3618.1904
$ws.Range("M40").Value = -3482.1904
$ws.Range("H46").Value = 3676.652
$ws.Range("I46").Value = 3916.5454
$ws.Range("K46").Value = 3916.5454
$ws.Range("M46").Value = -3728.5454
$ws.Range("H61").Value = 5525.1562
$ws.Range("J61").Value = 11404.454
$ws.Range("L61").Value = 11404.454
$ws.Range("N61").Value = -11808.454
$ws.Range("H93").Value = 2208.2856
$ws.Range("I93").Value = 909.6667
$ws.Range("K93").Value = 909.6667
$ws.Range("M93").Value = 338.3333
$ws.Range("H104").Value = 33437.8
$ws.Range("J104").Value = 33437.8
$ws.Range("L104").Value = 33437.8
$ws.Range("N104").Value = -40425.8
$ws.Range("H113").Value = 5525.1562
$ws.Range("J113").Value = 11404.454
$ws.Range("L113").Value = 11404.454
$ws.Range("N113").Value = -15744.454
$ws.Range("H122").Value = 4528
$ws.Range("I122").Value = 3833.074
$ws.Range("J122").Value = 6873.375
$ws.Range("K122").Value = 11499.222
$ws.Range("L122").Value = 20620.125
$ws.Range("M122").Value = -9049.222
$ws.Range("N122").Value = -25520.125
$ws.Range("H136").Value = 6674.853
$ws.Range("I136").Value = 3169.0688
$ws.Range("K136").Value = 9507.206399999999
$ws.Range("M136").Value = -6957.206399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 26999.8
$ws.Range("J4").Value = 38333.332
$ws.Range("L4").Value = 38333.332
$ws.Range("N4").Value = -38559.332
$ws.Range("H107").Value = 3406.0977
$ws.Range("I107").Value = 2764.7188
$ws.Range("J107").Value = 5686.5557
$ws.Range("K107").Value = 8294.1564
$ws.Range("L107").Value = 17059.6671
$ws.Range("M107").Value = -6374.1564
$ws.Range("N107").Value = -20899.6671
$ws.Range("H122").Value = 2706.5
$ws.Range("I122").Value = 2757.8333
$ws.Range("J122").Value = 2398.5
$ws.Range("K122").Value = 8273.499899999999
$ws.Range("L122").Value = 7195.5
$ws.Range("M122").Value = -5823.499899999999
$ws.Range("N122").Value = -12095.5
$ws.Range("H126").Value = 3727.05
$ws.Range("I126").Value = 3293.1428
$ws.Range("J126").Value = 4739.5
$ws.Range("K126").Value = 9879.428400000001
$ws.Range("L126").Value = 14218.5
$ws.Range("M126").Value = -7409.428400000001
$ws.Range("N126").Value = -19158.5
$ws.Range("H136").Value = 3448.4856
$ws.Range("I136").Value = 3160.1428
$ws.Range("J136").Value = 4601.857
$ws.Range("K136").Value = 9480.428400000001
$ws.Range("L136").Value = 13805.571
$ws.Range("M136").Value = -6930.428400000001
$ws.Range("N136").Value = -18905.571
